# "Actualizo IPC-AGo23 y opex 1er semestre 2023"
#
# Applies:
#  - IPC-Seriemensual: append August-2023 row (row 81)
#  - Aperturas: refresh regional opex % values (rows 2-10, cols B:H)
#  - IPC-DIC-Div: refresh regional IPC values (rows 2-14, cols B:H)
#  - IPC-Interanual: refresh regional IPC index values (rows 2-14, cols B:H)
#  - Active-tab moves from Impo-ICA to Aperturas

$wb = $excel.ActiveWorkbook

function Set-RowValues {
    param($ws, [int]$row, [object[]]$values)
    $n = $values.Length
    $arr = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) { $arr[0, $i] = $values[$i] }
    $startCol = 2  # column B
    $endCol = $startCol + $n - 1
    $rng = $ws.Range($ws.Cells.Item($row, $startCol), $ws.Cells.Item($row, $endCol))
    $rng.Value = $arr
}

# ---------------------------------------------------------------------------
# 1. IPC-Seriemensual: add the new monthly row (August 2023)
# ---------------------------------------------------------------------------
$wsSerie = $wb.Worksheets.Item("IPC-Seriemensual")
$wsSerie.Range("A81").Value = 45139
$wsSerie.Range("A81").NumberFormat = "dd/mm/yyyy"
$wsSerie.Range("B81").Value = 12.4
$wsSerie.Range("C81").Value = 10.7
$wsSerie.Range("D81").Value = 13.8
$wsSerie.Range("E81").Value = 8.3000000000000007

# Move the sheet's active selection to the newly entered row
$wsSerie.Range("B81:E81").Select()

# ---------------------------------------------------------------------------
# 2. Aperturas: updated regional opening percentages
# ---------------------------------------------------------------------------
$wsAperturas = $wb.Worksheets.Item("Aperturas")
Set-RowValues $wsAperturas 2 @(11.6, 12, 13, 14.7, 13.1, 12.5, 12.161660784025695)
Set-RowValues $wsAperturas 3 @(25.6, 23.8, 25, 26.6, 24.4, 18, 24.630684811291736)
Set-RowValues $wsAperturas 4 @(8.3000000000000007, 8.3000000000000007, 9.9, 10.8, 7.8, 9.8000000000000007, 8.590329761714921)
Set-RowValues $wsAperturas 5 @(7.8, 8.3000000000000007, 9, 7.2, 7.3, 9.6, 8.054183204113663)
Set-RowValues $wsAperturas 6 @(11.1, 11.6, 8, 13.1, 10.7, 10, 11.211081299883951)
Set-RowValues $wsAperturas 7 @(19.3, 17.3, 20.9, 16.7, 16.899999999999999, 15, 18.206750097855597)
Set-RowValues $wsAperturas 8 @(9.6, 12.3, 18.399999999999999, 16, 12.1, 11.4, 11.789156699560422)
Set-RowValues $wsAperturas 9 @(12.6, 11.4, 14.8, 11.9, 10.199999999999999, 14.5, 12.17504677347241)
Set-RowValues $wsAperturas 10 @(4.5, 10, 13.2, 11.7, 10.1, 9.1999999999999993, 7.9426590672652786)

# ---------------------------------------------------------------------------
# 3. IPC-DIC-Div: updated regional IPC (monthly variation) values
# ---------------------------------------------------------------------------
$wsDicDiv = $wb.Worksheets.Item("IPC-DIC-Div")
Set-RowValues $wsDicDiv 2 @(12.4, 12.3, 12.2, 14.2, 13.7, 12.3, 12.1)
Set-RowValues $wsDicDiv 3 @(15.6, 15.2, 15.6, 17.600000000000001, 17.899999999999999, 15.9, 13.8)
Set-RowValues $wsDicDiv 4 @(8.5, 7.4, 9, 11.8, 10.3, 8.6999999999999993, 8.1)
Set-RowValues $wsDicDiv 5 @(9.1, 8.3000000000000007, 9.9, 11.7, 9.9, 8.6, 9.6999999999999993)
Set-RowValues $wsDicDiv 6 @(9.1, 11.8, 7.1, 7.1, 5.0999999999999996, 7, 7.5)
Set-RowValues $wsDicDiv 7 @(14.1, 14.4, 13.4, 16.2, 15.6, 11.1, 14.5)
Set-RowValues $wsDicDiv 8 @(15.3, 16.100000000000001, 14.5, 17.899999999999999, 14.3, 13.2, 14.6)
Set-RowValues $wsDicDiv 9 @(10.5, 9.3000000000000007, 11.2, 11, 11.3, 11.2, 12.4)
Set-RowValues $wsDicDiv 10 @(4.5, 5, 3.7, 6.3, 4.4000000000000004, 3.6, 4.5999999999999996)
Set-RowValues $wsDicDiv 11 @(11.6, 11.5, 10.8, 12.9, 10.8, 13, 16.899999999999999)
Set-RowValues $wsDicDiv 12 @(8.6999999999999993, 9.5, 7, 5.4, 10.8, 8.6, 10.9)
Set-RowValues $wsDicDiv 13 @(12.4, 12.7, 11.4, 12.4, 15.1, 14.7, 11.1)
Set-RowValues $wsDicDiv 14 @(9.4, 9.1999999999999993, 9.1, 12.7, 9.3000000000000007, 9.4, 9.1999999999999993)

# ---------------------------------------------------------------------------
# 4. IPC-Interanual: updated regional IPC (year-over-year index) values
# ---------------------------------------------------------------------------
$wsInteranual = $wb.Worksheets.Item("IPC-Interanual")
Set-RowValues $wsInteranual 2 @(124.4, 125.2, 123.6, 126.1, 126, 119.7, 122.5)
Set-RowValues $wsInteranual 3 @(133.5, 135.6, 131.80000000000001, 131.4, 135.30000000000001, 130.1, 129.19999999999999)
Set-RowValues $wsInteranual 4 @(123.6, 118.4, 127.3, 140.69999999999999, 129.30000000000001, 121.5, 122.5)
Set-RowValues $wsInteranual 5 @(108.9, 110.3, 108.1, 111, 108.4, 98.7, 111.6)
Set-RowValues $wsInteranual 6 @(119.2, 128, 109.8, 133.5, 109.7, 105.9, 111.2)
Set-RowValues $wsInteranual 7 @(125, 124.8, 122.7, 134.9, 131.30000000000001, 120, 128.9)
Set-RowValues $wsInteranual 8 @(127.8, 124.6, 129, 135.5, 133, 131.6, 127.1)
Set-RowValues $wsInteranual 9 @(106.9, 106.5, 107.9, 102.9, 107, 102.4, 112.5)
Set-RowValues $wsInteranual 10 @(120.2, 122.6, 121.1, 115.2, 118.3, 111.5, 106.6)
Set-RowValues $wsInteranual 11 @(125.1, 125.9, 123.3, 125.9, 115.9, 129.1, 141.6)
Set-RowValues $wsInteranual 12 @(118.1, 125.5, 112.7, 99.4, 117.2, 106.1, 106.6)
Set-RowValues $wsInteranual 13 @(142.19999999999999, 136.6, 149.30000000000001, 138.9, 147.9, 139.9, 141.1)
Set-RowValues $wsInteranual 14 @(117, 119.6, 113, 124.4, 115, 116.3, 117)

# ---------------------------------------------------------------------------
# 5. Move the active tab from Impo-ICA to Aperturas (last sheet touched/
#    activated ends up as the workbook's active sheet on save)
# ---------------------------------------------------------------------------
$wsAperturas.Activate()
